# Generate Report for Handoff
# Adds a new handoff entry (ffff84959b15-56c7-4d03-8fb6-fe85b2821b63.md) alongside the
# existing one (renamed 2f8cc5e9-... -> d0da193e-...), pushing the
# ".localization-config" / "Not to be localized" row down by one on every sheet.

$wb = $excel.ActiveWorkbook

$oldMd = "2f8cc5e9-36bd-4329-8fa0-75d10cbea17f.md"
$newMd = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.md"
$newMd2 = "ffff84959b15-56c7-4d03-8fb6-fe85b2821b63.md"

$oldHashZh = "2f8cc5e9-36bd-4329-8fa0-75d10cbea17f.21f0e2e366d655d60e21386f30c24a2104a48801.zh-cn.xlf"
$newHashZh = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.zh-cn.xlf"
$oldHashDe = "2f8cc5e9-36bd-4329-8fa0-75d10cbea17f.21f0e2e366d655d60e21386f30c24a2104a48801.de-de.xlf"
$newHashDe = "d0da193e-f564-4aa1-9caa-cd92e2e17c5e.8614f8fdde6b09b065ea7020fa16dfa3020e6125.de-de.xlf"

$oldTimeZh = "2016-03-09 10:03:56"
$newTimeZh = "2016-03-09 10:05:07"
$oldTimeDe = "2016-03-09 10:04:05"
$newTimeDe = "2016-03-09 10:05:19"

$readyForHandoff = "Ready for handoff"
$notLocalized = "Not to be localized"
$configName = ".localization-config"
$include = "Include"
$ignored = "Ignored"
$epoch = "0001-01-01 00:00:00"

$baseRepo = "https://github.com/OpenLocalizationTest/oltest/blob/07c411c2591d0fcc359bbf713acd6fb33515a027"
$handoffZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/584437710a12a525dba2fa965818bac5de50e5ec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$handoffDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ee85d190c8b0e4279c1e479d5f045780658d772/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newMd
$ws1.Range("B2").Value = $readyForHandoff
$ws1.Range("C2").Value = $readyForHandoff

$ws1.Range("A3").Value = $newMd2
$ws1.Range("B3").Value = $readyForHandoff
$ws1.Range("C3").Value = $readyForHandoff

$ws1.Range("A4").Value = $configName
$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), ($baseRepo + "/e2e/" + $newMd), "", "", $newMd) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($baseRepo + "/e2e/" + $newMd2), "", "", $newMd2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($baseRepo + "/" + $configName), "", "", $configName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMd
$ws2.Range("B2").Value = $readyForHandoff
$ws2.Range("C2").Value = $newHashZh
$ws2.Range("D2").Value = $newTimeZh
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = $include

$ws2.Range("A3").Value = $newMd2
$ws2.Range("B3").Value = $readyForHandoff
$ws2.Range("C3").Value = $newHashZh
$ws2.Range("D3").Value = $newTimeZh
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $include

$ws2.Range("A4").Value = $configName
$ws2.Range("B4").Value = $notLocalized
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = $ignored

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), ($baseRepo + "/e2e/" + $newMd), "", "", $newMd) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($handoffZhBase + "/" + $newHashZh), "", "", $newHashZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($baseRepo + "/e2e/" + $newMd2), "", "", $newMd2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), ($handoffZhBase + "/" + $newHashZh), "", "", $newHashZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), ($baseRepo + "/" + $configName), "", "", $configName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMd
$ws3.Range("B2").Value = $readyForHandoff
$ws3.Range("C2").Value = $newHashDe
$ws3.Range("D2").Value = $newTimeDe
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = $include

$ws3.Range("A3").Value = $newMd2
$ws3.Range("B3").Value = $readyForHandoff
$ws3.Range("C3").Value = $newHashDe
$ws3.Range("D3").Value = $newTimeDe
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $include

$ws3.Range("A4").Value = $configName
$ws3.Range("B4").Value = $notLocalized
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = $ignored

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), ($baseRepo + "/e2e/" + $newMd), "", "", $newMd) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($handoffDeBase + "/" + $newHashDe), "", "", $newHashDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($baseRepo + "/e2e/" + $newMd2), "", "", $newMd2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), ($handoffDeBase + "/" + $newHashDe), "", "", $newHashDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), ($baseRepo + "/" + $configName), "", "", $configName) | Out-Null

Write-Output "Handback report rows generated."
